# Rename ObjTables document/table header attributes from
# PascalCase/UpperCamelCase to lowerCamelCase, per commit:
#   "changing document, table attributes to lowerCamelCase"
#
# Only the textual header cells (row 1 on the "!!Test" sheet, plus the
# "!!ObjTables Type=... Id=..." header row on every data sheet) carry the
# attribute names that need to change - everything else (styles, layout,
# other data) is left untouched.

$wb = $excel.ActiveWorkbook

$wsTest = $wb.Worksheets.Item("!!Test")
$wsTest.Range("A1").Value = "!!!ObjTables objTablesVersion='0.0.8'"
$wsTest.Range("A2").Value = "!!ObjTables type='Data' id='Test'"

$wsProperty = $wb.Worksheets.Item("!!Property")
$wsProperty.Range("A1").Value = "!!ObjTables type='Data' id='Property'"

$wsSubtests = $wb.Worksheets.Item("!!Subtests")
$wsSubtests.Range("A1").Value = "!!ObjTables type='Data' id='Subtest'"

$wsReferences = $wb.Worksheets.Item("!!References")
$wsReferences.Range("A1").Value = "!!ObjTables type='Data' id='Reference'"
